# Loan RBI, Variable Instalments
#
# The "Repayment Schedule" sheet gets a new (blank) column inserted right
# before the existing "Late" column, pushing "Late" / "Heading" /
# "Outstanding" one column to the right (N -> O, O -> P, P -> Q). The user
# then clicks over to the "Repayment Schedule" tab (making it the active
# sheet instead of "Input") and leaves the selection on cell S8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("Repayment Schedule")

# Insert a new blank column at N, shifting "Late"/"Heading"/"Outstanding"
# (and all the data under them) one column to the right.
$ws.Columns("N").Insert()

# Switch to the "Repayment Schedule" tab (it becomes the active sheet;
# "Input" - previously active - loses that status automatically) and
# leave the selection where the user ended up.
$ws.Activate() | Out-Null
$ws.Range("S8").Select() | Out-Null
